$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 1 ("Estrutura Analítica" title slide) - TextBox 3 (authors)
# Split " Gomes Júnior," into " " + "Gomes " + "Junior" and merge the
# leading comma into the following "Lucas Carvalho ..." run.
# -----------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$authorsShape = $s1.Shapes.Item(3)
$authorsRange = $authorsShape.TextFrame.TextRange

$oldRun = $authorsRange.Find(" Gomes Júnior,")
$oldRun.Text = " "

$spacePos = $authorsRange.Characters($oldRun.Start, 1)
$spacePos.InsertAfter("Gomes Junior")

$gomesRange = $authorsRange.Find("Gomes ")
$gomesRange.Text = $gomesRange.Text
$juniorRange = $authorsRange.Find("Junior")
$juniorRange.Text = $juniorRange.Text

$trailingSpace = $authorsRange.Find(" ", $juniorRange.Start)
$trailingSpace.Text = ""
$lucasRange = $authorsRange.Find("Lucas Carvalho Ribeiro, Pedro Henrique Gasparetto Lugão")
$lucasRange.Text = ", Lucas Carvalho Ribeiro, Pedro Henrique Gasparetto Lugão"

# -----------------------------------------------------------------
# Slide 2 ("Resumo do projeto") - Content Placeholder 2
# Split the last bullet's final word into its own run and append a
# new paragraph describing the tools used during development.
# -----------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item(2)
$contentRange = $contentShape.TextFrame.TextRange

$incrementalRange = $contentRange.Find("incremental")
$incrementalRange.Text = $incrementalRange.Text

$newParaText = "Durante o processo, utilizaremos a ferramenta Unity3D para o desenvolvimento, Pencil para o eventuais desenhos de interfaces, StarUML para o desenvolvimento de diagramas UML e GitHub para o controle de versão. Para demais tarefas os recursos serão definidos nas fases de análise de viabilidade, assim como os membros do grupo responsáveis por cada tarefa"
$contentRange.InsertAfter("`r" + $newParaText)

$newPara = $contentRange.Paragraphs(5, 1)
$paraStart = $newPara.Start

$pencilRange = $contentRange.Find("Pencil", $paraStart)
$pencilRange.Text = $pencilRange.Text

$starUmlRange = $contentRange.Find("StarUML", $paraStart)
$starUmlRange.Text = $starUmlRange.Text

$gitHubRange = $contentRange.Find("GitHub", $paraStart)
$gitHubRange.Text = $gitHubRange.Text

# Mark the placeholder to shrink text on overflow (closest reachable
# approximation of the authored <a:normAutofit fontScale=".." .../>).
$contentShape.TextFrame.AutoSize = 2
